$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item(1)

# Add the new sheet after the existing "AddPerson_TestData" sheet.
$ws2 = $wb.Worksheets.Add($null, $ws1)
$ws2.Name = "AddNewCase_TestData"

# Header row
$ws2.Range("A1").Value = "CaseRelatesTo"
$ws2.Range("B1").Value = "Name"
$ws2.Range("C1").Value = "Description"
$ws2.Range("D1").Value = "Tags"
$ws2.Range("E1").Value = "Track"

# Data row
$ws2.Range("A2").Value = "QA Test"
$ws2.Range("B2").Value = "Test Lost Car"
$ws2.Range("C2").Value = "Test Automation"
$ws2.Range("D2").Value = "Car"
$ws2.Range("E2").Value = " "

# Selection on the first sheet moves to B2, no longer the active tab.
$ws1.Range("B2").Select()

# New sheet becomes active, selection on E24.
$ws2.Range("E24").Select()
$ws2.Activate()
